$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the header row cells: *_old -> *_FV2310, *_new -> *_FV2404 ---
$headersFV2310 = @(
    "Segmentname_FV2310",
    "Segmentgruppe_FV2310",
    "Segment_FV2310",
    "Datenelement_FV2310",
    "Segment ID_FV2310",
    "Code_FV2310",
    "Qualifier_FV2310",
    "Beschreibung_FV2310",
    "Bedingungsausdruck_FV2310",
    "Bedingung_FV2310"
)
for ($i = 0; $i -lt $headersFV2310.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headersFV2310[$i]
}

$ws.Cells.Item(1, 11).Value = "diff"

$headersFV2404 = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)
for ($i = 0; $i -lt $headersFV2404.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $headersFV2404[$i]
}

# --- Turn the used range A1:U93 into an Excel Table (ListObject) ---
$tableRange = $ws.Range("A1:U93")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- Freeze the header row (pane split after row 1) ---
$ws.Range("A2").Select()
$win = $excel.ActiveWindow
$win.FreezePanes = $true
